# Add the new "banner_image" row to the Database sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

$ws.Range("A18").Value = "banner_image"
$ws.Range("B18").Value = "field for motobikes"

# Match the author's resulting view/selection state.
$ws.Activate()
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 7
